# Remove the three account rows that were dropped from the "Export" sheet
# (accounts 004547722/MARCIA, 004971448/CLOVIS, 004479965/DIEGO), shifting
# the remaining rows up - matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid as we go.
$ws.Rows.Item(6).Delete()   # 004479965 - DIEGO   - 3003.18
$ws.Rows.Item(5).Delete()   # 004971448 - CLOVIS   - 13000
$ws.Rows.Item(2).Delete()   # 004547722 - MARCIA   - 25872.83
